# Auto-generated: apply scheduled-runner market price updates to Shinryu_Profits sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 8080259
$ws.Range("I19").Value = 8537835
$ws.Range("J19").Value = 7693078.5
$ws.Range("K19").Value = 8537835
$ws.Range("L19").Value = 7693078.5
$ws.Range("M19").Value = -8537660
$ws.Range("N19").Value = -7693428.5
$ws.Range("H40").Value = 2800
$ws.Range("I40").Value = 2712.5
$ws.Range("J40").Value = 2916.6667
$ws.Range("K40").Value = 2712.5
$ws.Range("L40").Value = 2916.6667
$ws.Range("M40").Value = -2537.5
$ws.Range("N40").Value = -3266.6667
$ws.Range("H62").Value = 33336178
$ws.Range("J62").Value = 2859
$ws.Range("L62").Value = 2859
$ws.Range("N62").Value = -4107
$ws.Range("H65").Value = 33336178
$ws.Range("J65").Value = 2859
$ws.Range("L65").Value = 14295
$ws.Range("N65").Value = -20535
$ws.Range("H125").Value = 918.8570999999999
$ws.Range("I125").Value = 760.6667
$ws.Range("K125").Value = 6846.0003
$ws.Range("M125").Value = -4386.0003
$ws.Range("H127").Value = 888.59576
$ws.Range("J127").Value = 1028.8379
$ws.Range("L127").Value = 3086.5137
$ws.Range("N127").Value = -13006.5137
$ws.Range("H129").Value = 985.14
$ws.Range("I129").Value = 357.81818
$ws.Range("J129").Value = 1062.6742
$ws.Range("K129").Value = 1073.45454
$ws.Range("L129").Value = 3188.0226
$ws.Range("M129").Value = 3926.54546
$ws.Range("N129").Value = -13188.0226
$ws.Range("H131").Value = 23856.355
$ws.Range("I131").Value = 30438.295
$ws.Range("K131").Value = 91314.88499999999
$ws.Range("M131").Value = -86274.88499999999
$ws.Range("H132").Value = 1577.0571
$ws.Range("I132").Value = 1579.3235
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 4737.970499999999
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -2207.970499999999
$ws.Range("N132").Value = -9560
$ws.Range("H135").Value = 1266.6666
$ws.Range("I135").Value = 1266.6666
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 11399.9994
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -8864.999400000001
$ws.Range("N135").ClearContents()
$ws.Range("H137").Value = 10418754
$ws.Range("I137").Value = 20834956
$ws.Range("K137").Value = 62504868
$ws.Range("M137").Value = -62502318
$ws.Range("H138").Value = 4137.982
$ws.Range("I138").Value = 1312
$ws.Range("K138").Value = 3936
$ws.Range("M138").Value = 1204
$ws.Range("H139").Value = 11552.883
$ws.Range("J139").Value = 11552.883
$ws.Range("L139").Value = 11552.883
$ws.Range("N139").Value = -21832.883

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3430.4736
$ws.Range("I61").Value = 1054.6428
$ws.Range("K61").Value = 1054.6428
$ws.Range("M61").Value = -842.6428000000001
$ws.Range("H74").Value = 7306.3335
$ws.Range("I74").Value = 7468.75
$ws.Range("K74").Value = 7468.75
$ws.Range("M74").Value = -6594.75
$ws.Range("H77").Value = 7306.3335
$ws.Range("I77").Value = 7468.75
$ws.Range("K77").Value = 37343.75
$ws.Range("M77").Value = -32975.75
$ws.Range("H110").Value = 947.8570999999999
$ws.Range("I110").Value = 725.5
$ws.Range("J110").Value = 1763.1666
$ws.Range("K110").Value = 725.5
$ws.Range("L110").Value = 1763.1666
$ws.Range("M110").Value = 1319.5
$ws.Range("N110").Value = -5853.1666
$ws.Range("H122").Value = 1459.0588
$ws.Range("I122").Value = 1362.75
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 4088.25
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -1638.25
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 2843.5
$ws.Range("I132").Value = 2226.5
$ws.Range("K132").Value = 6679.5
$ws.Range("M132").Value = -4149.5
$ws.Range("H136").Value = 3430.4736
$ws.Range("I136").Value = 1054.6428
$ws.Range("K136").Value = 3163.9284
$ws.Range("M136").Value = -613.9284000000002

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 6973.5713
$ws.Range("J86").Value = 3561.6
$ws.Range("L86").Value = 3561.6
$ws.Range("N86").Value = -5807.6
$ws.Range("H89").Value = 6973.5713
$ws.Range("J89").Value = 3561.6
$ws.Range("L89").Value = 17808
$ws.Range("N89").Value = -29040
$ws.Range("H107").Value = 679.0303
$ws.Range("I107").Value = 487.68182
$ws.Range("J107").Value = 1061.7273
$ws.Range("K107").Value = 487.68182
$ws.Range("L107").Value = 1061.7273
$ws.Range("M107").Value = 1432.31818
$ws.Range("N107").Value = -4901.7273
$ws.Range("H132").Value = 3148.9
$ws.Range("I132").Value = 1946.1666
$ws.Range("J132").Value = 4953
$ws.Range("K132").Value = 5838.4998
$ws.Range("L132").Value = 14859
$ws.Range("M132").Value = -3308.4998
$ws.Range("N132").Value = -19919
$ws.Range("H134").Value = 2265.6333
$ws.Range("I134").Value = 1389.25
$ws.Range("J134").Value = 5771.1665
$ws.Range("K134").Value = 4167.75
$ws.Range("L134").Value = 17313.4995
$ws.Range("M134").Value = -1632.75
$ws.Range("N134").Value = -22383.4995

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 326.75
$ws.Range("I23").Value = 319.83334
$ws.Range("J23").Value = 333.66666
$ws.Range("K23").Value = 959.5000200000001
$ws.Range("L23").Value = 1000.99998
$ws.Range("M23").Value = -724.5000200000001
$ws.Range("N23").Value = -1470.99998
$ws.Range("H97").Value = 510.45
$ws.Range("I97").Value = 466.42856
$ws.Range("J97").Value = 534.1539
$ws.Range("K97").Value = 1399.28568
$ws.Range("L97").Value = 1602.4617
$ws.Range("M97").Value = -903.28568
$ws.Range("N97").Value = -2594.4617

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 177.43478
$ws.Range("I107").Value = 190.94737
$ws.Range("J107").Value = 113.25
$ws.Range("K107").Value = 190.94737
$ws.Range("L107").Value = 113.25
$ws.Range("M107").Value = 1729.05263
$ws.Range("N107").Value = -3953.25
$ws.Range("H113").Value = 13312.223
$ws.Range("I113").Value = 1279.6
$ws.Range("J113").Value = 28353
$ws.Range("K113").Value = 1279.6
$ws.Range("L113").Value = 28353
$ws.Range("M113").Value = 890.4000000000001
$ws.Range("N113").Value = -32693
$ws.Range("H126").Value = 2948.7896
$ws.Range("I126").Value = 2662.8462
$ws.Range("J126").Value = 3568.3333
$ws.Range("K126").Value = 7988.5386
$ws.Range("L126").Value = 10704.9999
$ws.Range("M126").Value = -5518.5386
$ws.Range("N126").Value = -15644.9999
$ws.Range("H132").Value = 4171.72
$ws.Range("I132").Value = 3927.111
$ws.Range("J132").Value = 4800.7144
$ws.Range("K132").Value = 11781.333
$ws.Range("L132").Value = 14402.1432
$ws.Range("M132").Value = -9251.332999999999
$ws.Range("N132").Value = -19462.1432

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 30000
$ws.Range("J59").Value = 30000
$ws.Range("L59").Value = 30000
$ws.Range("N59").Value = -31308
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2451.5
$ws.Range("I96").Value = 1760
$ws.Range("K96").Value = 1760
$ws.Range("M96").Value = -387
$ws.Range("H113").Value = 554.25
$ws.Range("I113").Value = 490.6
$ws.Range("K113").Value = 1471.8
$ws.Range("M113").Value = 698.1999999999998
